# Normalize district-name casing on the "İş Takip Listesi" sheet (rows 95-122,
# column B) and on the "Güncelleme" sheet (rows 2-29, column A) so they match
# the already-existing lowercase district strings used elsewhere in the
# workbook ("Akdeniz", "Toroslar", "Tarsus", "Mezitli", "Yenişehir" instead of
# the ALL-CAPS duplicates). This also updates the sheet views / selections and
# the ignoredErrors range on "Güncelleme" to match what Excel recomputes.

$wb = $excel.ActiveWorkbook

$wsTakip = $wb.Worksheets.Item("İş Takip Listesi")
$wsGuncelleme = $wb.Worksheets.Item("Güncelleme")

# Row (95-122) -> district name, taken from the target workbook state.
$districtByRow = @{
    95  = "Akdeniz"
    96  = "Akdeniz"
    97  = "Akdeniz"
    98  = "Toroslar"
    99  = "Toroslar"
    100 = "Toroslar"
    101 = "Toroslar"
    102 = "Toroslar"
    103 = "Toroslar"
    104 = "Toroslar"
    105 = "Tarsus"
    106 = "Tarsus"
    107 = "Tarsus"
    108 = "Tarsus"
    109 = "Tarsus"
    110 = "Tarsus"
    111 = "Tarsus"
    112 = "Tarsus"
    113 = "Tarsus"
    114 = "Tarsus"
    115 = "Tarsus"
    116 = "Mezitli"
    117 = "Toroslar"
    118 = "Toroslar"
    119 = "Toroslar"
    120 = "Toroslar"
    121 = "Toroslar"
    122 = "Yenişehir"
}

foreach ($row in 95..122) {
    $wsTakip.Range("B$row").Value = $districtByRow[$row]
}

# "Güncelleme" column A mirrors "İş Takip Listesi" column B for the same
# records (row N on Güncelleme corresponds to row N+93 on İş Takip Listesi).
foreach ($row in 2..29) {
    $wsGuncelleme.Range("A$row").Value = $districtByRow[$row + 93]
}

# Sheet view / selection updates captured by the diff.
$wsTakip.Application.ActiveWindow.ScrollRow = 100
$wsTakip.Range("A100").Select()
$wsTakip.Range("B95:B122").Select()

$wsGuncelleme.Activate()
$wsGuncelleme.Application.ActiveWindow.ScrollRow = 1
$wsGuncelleme.Range("A2:A29").Select()
